$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E12").Value = 1361974149.3
$ws.Range("F12").Value = 1303578329

$ws.Range("E13").Value = 332114255.61000007
$ws.Range("F13").Value = 325268233.60000002

$ws.Range("E14").Value = 537985.22000000253

$ws.Range("F16").Value = -53616441.740000002

$ws.Range("E18").Formula = "=SUM(E12:E17)"

$ws.Range("E19").Value = -384700000.00000012

$ws.Range("E21").Formula = "=SUM(E18:E20)"

$ws.Range("E26").Value = 1018613403.9980445
$ws.Range("F26").Value = 1026703455
